$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 279, shifting existing rows 279+ down to 281+
$ws.Rows("279:280").Insert()

# Populate the two newly inserted rows with fresh data
$ws.Range("A279").Value = 1
$ws.Range("B279").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C279").Value = "Arica y Parinacota"
$ws.Range("D279").Value = 44722
$ws.Range("E279").Value = 15
$ws.Range("F279").Value = 100112032
$ws.Range("G279").Value = "Zapallo italiano"
$ws.Range("H279").Value = "Huracán"
$ws.Range("I279").Value = "Primera"
$ws.Range("J279").Value = 160
$ws.Range("K279").Value = 6500
$ws.Range("L279").Value = 7000
$ws.Range("M279").Value = 6750
$ws.Range("N279").Value = "`$/caja 70 unidades"
$ws.Range("O279").Value = "Región de Arica y Parinacota"
$ws.Range("P279").Value = 96
$ws.Range("Q279").Value = 70
$ws.Range("R279").Value = "Hortaliza"

$ws.Range("A280").Value = 1
$ws.Range("B280").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C280").Value = "Arica y Parinacota"
$ws.Range("D280").Value = 44722
$ws.Range("E280").Value = 15
$ws.Range("F280").Value = 100112032
$ws.Range("G280").Value = "Zapallo italiano"
$ws.Range("H280").Value = "Huracán"
$ws.Range("I280").Value = "Segunda"
$ws.Range("J280").Value = 160
$ws.Range("K280").Value = 5500
$ws.Range("L280").Value = 6000
$ws.Range("M280").Value = 5750
$ws.Range("N280").Value = "`$/caja 100 unidades"
$ws.Range("O280").Value = "Región de Arica y Parinacota"
$ws.Range("P280").Value = 58
$ws.Range("Q280").Value = 100
$ws.Range("R280").Value = "Hortaliza"
